$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.07"
$ws.Range("E2").Value = "'-0.79%"
$ws.Range("D3").Value = "'27.32"
$ws.Range("E3").Value = "'1.84%"
$ws.Range("D4").Value = "'4.797"
$ws.Range("E4").Value = "'-2.39%"
$ws.Range("E5").Value = "'-0.60%"
$ws.Range("D6").Value = "'6.954"
$ws.Range("E6").Value = "'-0.16%"
$ws.Range("D7").Value = "'1.315"
$ws.Range("E7").Value = "'8.84%"
$ws.Range("D8").Value = "'0.8766"
$ws.Range("E8").Value = "'-1.06%"
$ws.Range("D9").Value = "'0.1543"
$ws.Range("E9").Value = "'4.07%"
$ws.Range("E10").Value = "'-2.89%"
$ws.Range("D11").Value = "'0.07586"
$ws.Range("E11").Value = "'2.38%"
$ws.Range("D12").Value = "'0.02983"
$ws.Range("E12").Value = "'-4.83%"
$ws.Range("E13").Value = "'-0.38%"
$ws.Range("D14").Value = "'0.001574"
$ws.Range("E14").Value = "'0.41%"
$ws.Range("D15").Value = "'0.0006415"
$ws.Range("E15").Value = "'1.29%"
$ws.Range("D16").Value = "'0.005792"
$ws.Range("E16").Value = "'-3.71%"
$ws.Range("E17").Value = "'-1.04%"
$ws.Range("D18").Value = "'3.303"
$ws.Range("E18").Value = "'-1.56%"
$ws.Range("D19").Value = "'2.272"
$ws.Range("E19").Value = "'-0.35%"
$ws.Range("E21").Value = "'1.74%"
$ws.Range("D22").Value = "'3.941"
$ws.Range("E22").Value = "'0.68%"
$ws.Range("D23").Value = "'0.04406"
$ws.Range("E23").Value = "'1.29%"
$ws.Range("D24").Value = "'0.001173"
$ws.Range("E24").Value = "'-0.48%"
$ws.Range("D25").Value = "'0.003864"
$ws.Range("E25").Value = "'5.09%"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("E27").Value = "'19.80%"
$ws.Range("D40").Value = "'0.04171"
$ws.Range("E40").Value = "'2.48%"
$ws.Range("D41").Value = "'0.006812"
$ws.Range("E41").Value = "'2.58%"
$ws.Range("E42").Value = "'0.44%"
$ws.Range("D43").Value = "'0.002018"
$ws.Range("E43").Value = "'-14.50%"
$ws.Range("D44").Value = "'0.01117"
$ws.Range("E44").Value = "'-12.99%"
$ws.Range("D45").Value = "'0.00005170"
$ws.Range("E45").Value = "'-1.74%"
$ws.Range("D46").Value = "'1.486"
$ws.Range("E46").Value = "'-36.89%"
$ws.Range("E47").Value = "'8.42%"
